$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" data column (P) to the right of the existing "2021" column (O),
# mirroring each row's existing O-column formatting onto the new P cell.
$ws.Range("O4:O14").Copy() | Out-Null
$ws.Range("P4:P14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("O16:O17").Copy() | Out-Null
$ws.Range("P16:P17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# Populate the new 2022 values.
$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 11.4
$ws.Range("P6").Value = 12.6
$ws.Range("P7").Value = 9.8000000000000007
$ws.Range("P8").Value = 11.4
$ws.Range("P9").Value = 5.4
$ws.Range("P10").Value = 4.7
$ws.Range("P11").Value = 3.4
$ws.Range("P12").Value = 17.7
$ws.Range("P13").Value = 20.5
$ws.Range("P14").Value = 8.4
$ws.Range("P16").Value = 12.9
$ws.Range("P17").Value = 10.5

# Match the saved selection from the edit.
$ws.Range("Q4").Select() | Out-Null
